# Se cambia el método para encontrar y seleccionar contacto:
# se agregan nuevos contactos/filas (4, 5 y 6) a la hoja "Hoja1" y se
# selecciona la celda B7 (próxima fila disponible) como en el flujo real
# de la macro que recorre los contactos.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fila 4 - Ferreteria 1
$ws.Range("A4").Value = "666"
$ws.Range("B4").Value = "Ferreteria 1"
$ws.Range("C4").Value = "sebas"
$ws.Range("D4").Value = "3227804602"
$ws.Range("E4").Value = "121212"
$ws.Range("F4").Value = "1111"

# Fila 5 - Ferreteria 2
$ws.Range("A5").Value = "555"
$ws.Range("B5").Value = "Ferreteria 2"
$ws.Range("C5").Value = "Melqui"
$ws.Range("F5").Value = "2222"
$ws.Range("D5").Value = "3176794454"
$ws.Range("E5").Value = "121212"

# Fila 6
$ws.Range("A6").Value = "1212"
$ws.Range("C6").Value = "asdasdfv"
$ws.Range("D6").Value = "4545454545"
$ws.Range("F6").Value = "3333"
$ws.Range("E6").Value = "121212"
$ws.Range("B6").Value = "imagine"

# Se selecciona la celda B7 (el método ahora busca y selecciona el
# siguiente contacto disponible en la columna B)
$ws.Range("B7").Select()
